# Apply the "Fixed bugs in Practice registry / Provider registry" data corrections:
# rows 2-6 get new CustID (A), NAME (B), and YEAR (C) values. The A and C columns
# hold digit-only values but must remain stored as text (matching the workbook's
# original inlineStr/text representation), so we prefix them with a leading
# apostrophe to force Excel to keep them as text instead of auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'1000"
$ws.Range("B2").Value = "Echelon Care Organisation"
$ws.Range("C2").Value = "'2022"

$ws.Range("A3").Value = "'200"
$ws.Range("B3").Value = "Enlighten Care Organisation"
$ws.Range("C3").Value = "'2022"

$ws.Range("A4").Value = "'1500"
$ws.Range("B4").Value = "Excalibur Physicians Association"
$ws.Range("C4").Value = "'2022"

$ws.Range("A5").Value = "'6500"
$ws.Range("B5").Value = "Guardian Care Delivery Organisation"
$ws.Range("C5").Value = "'2022"

$ws.Range("A6").Value = "'1100"
$ws.Range("B6").Value = "Saint Care Delivery Organisation"
$ws.Range("C6").Value = "'2022"
